# Update the "想去人数" (column F) counts on the 展览 and 全部类型 sheets
# to reflect the latest scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 897
$ws1.Range("F20").Value = 92
$ws1.Range("F22").Value = 3257
$ws1.Range("F23").Value = 5599
$ws1.Range("F25").Value = 5
$ws1.Range("F29").Value = 3202
$ws1.Range("F31").Value = 2403
$ws1.Range("F35").Value = 177
$ws1.Range("F40").Value = 871

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F16").Value = 897
$ws4.Range("F21").Value = 92
$ws4.Range("F23").Value = 3257
$ws4.Range("F24").Value = 5599
$ws4.Range("F26").Value = 5
$ws4.Range("F30").Value = 3202
$ws4.Range("F32").Value = 2403
$ws4.Range("F36").Value = 177
$ws4.Range("F41").Value = 871
